$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 153927
$ws.Range("C5").Value = 8607
$ws.Range("C6").Value = 341
$ws.Range("C7").Value = 5.59
